$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-11-13 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-14 Tuesday", 2)

$d.Content.Find.Execute("32×33=", $true, $false, $false, $false, $false, $true, 1, $false, "19×91=", 2)
$d.Content.Find.Execute("86×35=", $true, $false, $false, $false, $false, $true, 1, $false, "68×86=", 2)
$d.Content.Find.Execute("60×15=", $true, $false, $false, $false, $false, $true, 1, $false, "11×13=", 2)
$d.Content.Find.Execute("66×46=", $true, $false, $false, $false, $false, $true, 1, $false, "27×40=", 2)
$d.Content.Find.Execute("95×56=", $true, $false, $false, $false, $false, $true, 1, $false, "22×70=", 2)

$d.Content.Find.Execute("99×56=", $true, $false, $false, $false, $false, $true, 1, $false, "11×46=", 2)
$d.Content.Find.Execute("47×14=", $true, $false, $false, $false, $false, $true, 1, $false, "45×82=", 2)
$d.Content.Find.Execute("64×17=", $true, $false, $false, $false, $false, $true, 1, $false, "23×37=", 2)
$d.Content.Find.Execute("74×65=", $true, $false, $false, $false, $false, $true, 1, $false, "99×58=", 2)
$d.Content.Find.Execute("62×61=", $true, $false, $false, $false, $false, $true, 1, $false, "85×13=", 2)

$d.Content.Find.Execute("98×28=", $true, $false, $false, $false, $false, $true, 1, $false, "24×70=", 2)
$d.Content.Find.Execute("39×91=", $true, $false, $false, $false, $false, $true, 1, $false, "88×68=", 2)
$d.Content.Find.Execute("64×67=", $true, $false, $false, $false, $false, $true, 1, $false, "16×38=", 2)
$d.Content.Find.Execute("37×93=", $true, $false, $false, $false, $false, $true, 1, $false, "47×87=", 2)
$d.Content.Find.Execute("17×19=", $true, $false, $false, $false, $false, $true, 1, $false, "17×32=", 2)

$d.Content.Find.Execute("11×42=", $true, $false, $false, $false, $false, $true, 1, $false, "78×82=", 2)
$d.Content.Find.Execute("37×16=", $true, $false, $false, $false, $false, $true, 1, $false, "29×87=", 2)
$d.Content.Find.Execute("38×33=", $true, $false, $false, $false, $false, $true, 1, $false, "51×15=", 2)
$d.Content.Find.Execute("87×34=", $true, $false, $false, $false, $false, $true, 1, $false, "40×68=", 2)
$d.Content.Find.Execute("36×98=", $true, $false, $false, $false, $false, $true, 1, $false, "77×43=", 2)

$d.Content.Find.Execute("50×86=", $true, $false, $false, $false, $false, $true, 1, $false, "94×98=", 2)
$d.Content.Find.Execute("67×54=", $true, $false, $false, $false, $false, $true, 1, $false, "60×74=", 2)
$d.Content.Find.Execute("42×93=", $true, $false, $false, $false, $false, $true, 1, $false, "35×94=", 2)
$d.Content.Find.Execute("54×50=", $true, $false, $false, $false, $false, $true, 1, $false, "29×21=", 2)
$d.Content.Find.Execute("91×59=", $true, $false, $false, $false, $false, $true, 1, $false, "74×30=", 2)
